$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.140.90'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '1.826.44'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.84%  '
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4702'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3646'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07403'
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8800'
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.37'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '1.917.65'
$ws.Range('E12').Value = '  +5.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07335'
$ws.Range('E13').Value = '  +3.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.384'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.19'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.514'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008712'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.011'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').Value = '27.709.00'
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.65'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.244'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.58'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '2.105.77'
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.882'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.62'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.55'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.134'
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.171'
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.24'
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08935'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.166'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7408'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.508'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.941'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.011'
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05294'
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01945'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.407'
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.945'
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.186'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5259'
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1642'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.369'
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4878'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.38'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.011'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.40'
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.651'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06302'
$ws.Range('E51').Value = '  +0.07%  '
